# Update the "想去人数" (F column) counts across the four sheets to reflect
# a refreshed data pull, per the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(8, 6).Value = 1345
$ws.Cells.Item(9, 6).Value = 2032
$ws.Cells.Item(15, 6).Value = 3794
$ws.Cells.Item(18, 6).Value = 2860
$ws.Cells.Item(19, 6).Value = 746
$ws.Cells.Item(22, 6).Value = 68
$ws.Cells.Item(23, 6).Value = 1988
$ws.Cells.Item(27, 6).Value = 193
$ws.Cells.Item(28, 6).Value = 7993
$ws.Cells.Item(29, 6).Value = 5481
$ws.Cells.Item(32, 6).Value = 740
$ws.Cells.Item(40, 6).Value = 4577
$ws.Cells.Item(41, 6).Value = 804
$ws.Cells.Item(42, 6).Value = 52

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(18, 6).Value = 145

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 8161
$ws.Cells.Item(4, 6).Value = 1213

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 8161
$ws.Cells.Item(5, 6).Value = 1213
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(10, 6).Value = 1345
$ws.Cells.Item(14, 6).Value = 3794
$ws.Cells.Item(16, 6).Value = 2860
$ws.Cells.Item(17, 6).Value = 746
$ws.Cells.Item(20, 6).Value = 1988
$ws.Cells.Item(30, 6).Value = 193
$ws.Cells.Item(31, 6).Value = 7993
$ws.Cells.Item(32, 6).Value = 5481
$ws.Cells.Item(35, 6).Value = 740
$ws.Cells.Item(44, 6).Value = 4577
$ws.Cells.Item(45, 6).Value = 804
$ws.Cells.Item(46, 6).Value = 52
